# Add files via upload
# The schedule's week-of dates in columns A and B were re-stamped with an
# explicit year prefix ("5/23" -> "2025/5/23"), and a new trailing week
# ("2026/1/2") was appended to column B on the last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (rows 1-32): prepend "2025/" to the existing "<month>/<day>" text ---
for ($r = 1; $r -le 32; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cur = $cell.Value()
    if ($cur -ne $null -and $cur -ne "") {
        $cell.Value = "2025/" + $cur
    }
}

# --- Column A (rows 6-33): prepend "2025/" to the existing "<month>/<day>" text ---
for ($r = 6; $r -le 33; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cur = $cell.Value()
    if ($cur -ne $null -and $cur -ne "") {
        $cell.Value = "2025/" + $cur
    }
}

# --- New trailing week added to column B, row 33 ---
$ws.Cells.Item(33, 2).Value = "2026/1/2"

# --- B1:B3 previously used the wrap-text style (like column C); align them with
#     the plain text style already used by B4:B32 ---
$ws.Range("B4").Copy() | Out-Null
$ws.Range("B1:B3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Leave the selection where Excel would land after the edits ---
$ws.Range("C34").Select() | Out-Null
